$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1080.5
$ws.Range("I135").Value = 1047.4884
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 9427.3956
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -6892.3956
$ws.Range("N135").Value = -27570

$ws.Range("H137").Value = 83335430
$ws.Range("I137").Value = 200001470
$ws.Range("K137").Value = 600004410
$ws.Range("M137").Value = -600001860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 250985.5
$ws.Range("I28").Value = 250985.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 250985.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -250793.5
$ws.Range("N28").ClearContents()

$ws.Range("H61").Value = 2649.3225
$ws.Range("I61").Value = 1783.5
$ws.Range("K61").Value = 1783.5
$ws.Range("M61").Value = -1571.5

$ws.Range("H74").Value = 3944.2766
$ws.Range("I74").Value = 1133.7028
$ws.Range("J74").Value = 14343.4
$ws.Range("K74").Value = 1133.7028
$ws.Range("L74").Value = 14343.4
$ws.Range("M74").Value = -259.7028
$ws.Range("N74").Value = -16091.4

$ws.Range("H77").Value = 3944.2766
$ws.Range("I77").Value = 1133.7028
$ws.Range("J77").Value = 14343.4
$ws.Range("K77").Value = 5668.514
$ws.Range("L77").Value = 71717
$ws.Range("M77").Value = -1300.514
$ws.Range("N77").Value = -80453

$ws.Range("H96").Value = 30344
$ws.Range("J96").Value = 30344
$ws.Range("L96").Value = 30344
$ws.Range("N96").Value = -35836

$ws.Range("H99").Value = 250985.5
$ws.Range("I99").Value = 250985.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 250985.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -247990.5
$ws.Range("N99").ClearContents()

$ws.Range("H101").Value = 343734.66
$ws.Range("J101").Value = 343734.66
$ws.Range("L101").Value = 343734.66
$ws.Range("N101").Value = -350224.66

$ws.Range("H136").Value = 2649.3225
$ws.Range("I136").Value = 1783.5
$ws.Range("K136").Value = 5350.5
$ws.Range("M136").Value = -2800.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 726.4286
$ws.Range("J80").Value = 815.7059
$ws.Range("L80").Value = 815.7059
$ws.Range("N80").Value = -2811.7059

$ws.Range("H83").Value = 726.4286
$ws.Range("J83").Value = 815.7059
$ws.Range("L83").Value = 4078.5295
$ws.Range("N83").Value = -14062.5295

$ws.Range("H134").Value = 2991.8838
$ws.Range("I134").Value = 1994.5769
$ws.Range("J134").Value = 4517.1763
$ws.Range("K134").Value = 5983.7307
$ws.Range("L134").Value = 13551.5289
$ws.Range("M134").Value = -3448.7307
$ws.Range("N134").Value = -18621.5289

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1527.6487
$ws.Range("I31").Value = 907.8214
$ws.Range("J31").Value = 3456
$ws.Range("K31").Value = 907.8214
$ws.Range("L31").Value = 3456
$ws.Range("M31").Value = -612.8214
$ws.Range("N31").Value = -4046

$ws.Range("H34").Value = 1527.6487
$ws.Range("I34").Value = 907.8214
$ws.Range("J34").Value = 3456
$ws.Range("K34").Value = 907.8214
$ws.Range("L34").Value = 3456
$ws.Range("M34").Value = -705.8214
$ws.Range("N34").Value = -3860

$ws.Range("H58").Value = 1892.4048
$ws.Range("I58").Value = 827.88464
$ws.Range("J58").Value = 3622.25
$ws.Range("K58").Value = 827.88464
$ws.Range("L58").Value = 3622.25
$ws.Range("M58").Value = -624.88464
$ws.Range("N58").Value = -4028.25

$ws.Range("H99").Value = 8930485
$ws.Range("I99").Value = 20834660
$ws.Range("J99").Value = 2353.5
$ws.Range("K99").Value = 20834660
$ws.Range("L99").Value = 2353.5
$ws.Range("M99").Value = -20833162
$ws.Range("N99").Value = -5349.5

$ws.Range("H122").Value = 1824.3636
$ws.Range("I122").Value = 1081.1428
$ws.Range("K122").Value = 3243.4284
$ws.Range("M122").Value = -793.4284000000002

$ws.Range("H126").Value = 8930485
$ws.Range("I126").Value = 20834660
$ws.Range("J126").Value = 2353.5
$ws.Range("K126").Value = 62503980
$ws.Range("L126").Value = 7060.5
$ws.Range("M126").Value = -62501510
$ws.Range("N126").Value = -12000.5

$ws.Range("H136").Value = 1892.4048
$ws.Range("I136").Value = 827.88464
$ws.Range("J136").Value = 3622.25
$ws.Range("K136").Value = 2483.65392
$ws.Range("L136").Value = 10866.75
$ws.Range("M136").Value = 66.34608000000026
$ws.Range("N136").Value = -15966.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 254.42857
$ws.Range("I6").Value = 197.07692
$ws.Range("K6").Value = 591.23076
$ws.Range("M6").Value = -478.23076

$ws.Range("H122").Value = 630
$ws.Range("I122").Value = 272.25
$ws.Range("K122").Value = 2450.25
$ws.Range("M122").Value = -0.25

$ws.Range("H140").Value = 4065.5854
$ws.Range("I140").Value = 5096
$ws.Range("J140").Value = 2610.8823
$ws.Range("K140").Value = 15288
$ws.Range("L140").Value = 7832.646900000001
$ws.Range("M140").Value = -10108
$ws.Range("N140").Value = -18192.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8250
$ws.Range("I70").Value = 12700
$ws.Range("J70").Value = 5707.143
$ws.Range("K70").Value = 12700
$ws.Range("L70").Value = 5707.143
$ws.Range("M70").Value = -12430
$ws.Range("N70").Value = -6247.143

$ws.Range("H73").Value = 8250
$ws.Range("I73").Value = 12700
$ws.Range("J73").Value = 5707.143
$ws.Range("K73").Value = 12700
$ws.Range("L73").Value = 5707.143
$ws.Range("M73").Value = -11764
$ws.Range("N73").Value = -7579.143

$ws.Range("H102").Value = 4592.5713
$ws.Range("I102").Value = 3529.6
$ws.Range("J102").Value = 7250
$ws.Range("K102").Value = 3529.6
$ws.Range("L102").Value = 7250
$ws.Range("M102").Value = -1907.6
$ws.Range("N102").Value = -10494

$ws.Range("H104").Value = 150000
$ws.Range("J104").Value = 150000
$ws.Range("L104").Value = 150000
$ws.Range("N104").Value = -156988

$ws.Range("H126").Value = 2994.7646
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3060.7334
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 9182.200199999999
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -14122.2002

$ws.Range("H132").Value = 2640.3
$ws.Range("I132").Value = 2452.561
$ws.Range("J132").Value = 3495.5557
$ws.Range("K132").Value = 7357.683000000001
$ws.Range("L132").Value = 10486.6671
$ws.Range("M132").Value = -4827.683000000001
$ws.Range("N132").Value = -15546.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3689.2104
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws.Range("H136").Value = 4109.122
$ws.Range("I136").Value = 2678.8823
$ws.Range("J136").Value = 11056
$ws.Range("K136").Value = 8036.646900000001
$ws.Range("L136").Value = 33168
$ws.Range("M136").Value = -5486.646900000001
$ws.Range("N136").Value = -38268

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 202800.8
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 7938899.5
$ws.Range("I132").Value = 10419049
$ws.Range("J132").Value = 2419.5334
$ws.Range("K132").Value = 31257147
$ws.Range("L132").Value = 7258.600199999999
$ws.Range("M132").Value = -31254617
$ws.Range("N132").Value = -12318.6002

$ws.Range("H136").Value = 7776049
$ws.Range("I136").Value = 8359082.5
$ws.Range("J136").Value = 2263.3333
$ws.Range("K136").Value = 25077247.5
$ws.Range("L136").Value = 6789.999899999999
$ws.Range("M136").Value = -25074697.5
$ws.Range("N136").Value = -11889.9999
